$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet references
# ---------------------------------------------------------------------------
$tournaments = $wb.Worksheets.Item("tournaments")
$quests = $wb.Worksheets.Item("quests")

# ---------------------------------------------------------------------------
# 1) "tournaments" sheet - row 128: give it an explicit custom row height
# ---------------------------------------------------------------------------
$tournaments.Rows.Item(128).RowHeight = 14.25

# ---------------------------------------------------------------------------
# 2) "tournaments" sheet - rows 153/154/155 turn into three new TubeMan
#    tournament definitions (Kill Normal / Kill Time Limit / Kill Time Attack)
#
#    The shared-string table is append-only in write order, so the new
#    unique strings are entered in the exact order they appear in the
#    saved workbook: the three TID names first, then "TubeMan", then
#    "icon_tubeman" (the quest-sheet TID is entered later, see below).
# ---------------------------------------------------------------------------
$tournaments.Range("B153").Value = "TID_EVENT_TOURNAMENT_KILL_NORMAL_TUBE_MAN_FLOAT"
$tournaments.Range("B154").Value = "TID_EVENT_TOURNAMENT_KILL_TIME_LIMIT_TUBE_MAN_FLOAT"
$tournaments.Range("B155").Value = "TID_EVENT_TOURNAMENT_KILL_TIME_ATTACK_TUBE_MAN_FLOAT"

$tournaments.Range("E153").Value = "TubeMan"
$tournaments.Range("G153").Value = "icon_tubeman"

# -- Row 153 (Normal) --------------------------------------------------------
# Column A loses its old "<Definition>" marker/style and becomes a plain
# empty cell (style copied from A133, which already has that blank look).
$tournaments.Range("A133").Copy()
$tournaments.Range("A153").PasteSpecial(-4122)
$tournaments.Range("A153").Value = ""

$tournaments.Range("C153").Value = "kill"
$tournaments.Range("D153").Value = 0

# -- Row 154 (Time Limit) ----------------------------------------------------
$tournaments.Range("A154").Value = ""

# C154 and G154 need their style bumped to match the other rows in this block
$tournaments.Range("C150").Copy()
$tournaments.Range("C154").PasteSpecial(-4122)
$tournaments.Range("C154").Value = "kill"

$tournaments.Range("D154").Value = 2
$tournaments.Range("E154").Value = "TubeMan"

$tournaments.Range("G128").Copy()
$tournaments.Range("G154").PasteSpecial(-4122)
$tournaments.Range("G154").Value = "icon_tubeman"

# -- Row 155 (Time Attack) ---------------------------------------------------
$tournaments.Range("A155").Value = ""

$tournaments.Range("C150").Copy()
$tournaments.Range("C155").PasteSpecial(-4122)
$tournaments.Range("C155").Value = "kill"

$tournaments.Range("D155").Value = 1
$tournaments.Range("E155").Value = "TubeMan"

$tournaments.Range("G128").Copy()
$tournaments.Range("G155").PasteSpecial(-4122)
$tournaments.Range("G155").Value = "icon_tubeman"

# ---------------------------------------------------------------------------
# 3) "quests" sheet - add new row 54 (TubeMan global quest entry), formatted
#    like row 53 right above it.
# ---------------------------------------------------------------------------
$quests.Range("A53:F53").Copy()
$quests.Range("A54:F54").PasteSpecial(-4122)

$quests.Range("A54").Value = "<Definition>"
$quests.Range("B54").Value = "TID_GLOBAL_EVENT_EAT_TUBE_MAN_FLOAT"
$quests.Range("C54").Value = "kill"
$quests.Range("D54").Value = "TubeMan"
$quests.Range("F54").Value = "icon_tubeman"

# ---------------------------------------------------------------------------
# 4) View / selection state
#    - "tournaments": scroll down, select E153, no longer the active tab
#    - "quests": becomes the active tab, select F54
# ---------------------------------------------------------------------------
$tournaments.Activate()
$tournaments.Range("A122").Select()
$excel.ActiveWindow.ScrollRow = 122
$tournaments.Range("E153").Select()

$quests.Activate()
$excel.ActiveWindow.ScrollRow = 22
$quests.Range("F54").Select()
